$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3342
$ws.Range("I69").Value = 3342
$ws.Range("K69").Value = 10026
$ws.Range("M69").Value = -9152

$ws.Range("H72").Value = 3342
$ws.Range("I72").Value = 3342
$ws.Range("K72").Value = 30078
$ws.Range("M72").Value = -25710

$ws.Range("H107").Value = 419.57144
$ws.Range("I107").Value = 272.83334
$ws.Range("K107").Value = 272.83334
$ws.Range("M107").Value = 1647.16666

$ws.Range("H112").Value = 1471
$ws.Range("J112").Value = 1544.5
$ws.Range("L112").Value = 4633.5
$ws.Range("N112").Value = -6849.5

$ws.Range("H132").Value = 1319.5
$ws.Range("I132").Value = 1188.9445
$ws.Range("K132").Value = 3566.8335
$ws.Range("M132").Value = -1036.8335

$ws.Range("H137").Value = 3333.2104
$ws.Range("I137").Value = 934.4286
$ws.Range("J137").Value = 4732.5
$ws.Range("K137").Value = 2803.2858
$ws.Range("L137").Value = 14197.5
$ws.Range("M137").Value = -253.2857999999997
$ws.Range("N137").Value = -19297.5

$ws.Range("H138").Value = 3795.848
$ws.Range("J138").Value = 4336.9443
$ws.Range("L138").Value = 13010.8329
$ws.Range("N138").Value = -23290.8329

$ws.Range("H141").Value = 3347.7
$ws.Range("I141").Value = 3164.111
$ws.Range("K141").Value = 9492.332999999999
$ws.Range("M141").Value = -4312.332999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1147.9231
$ws.Range("I2").Value = 993.9091
$ws.Range("K2").Value = 993.9091
$ws.Range("M2").Value = -880.9091

$ws.Range("H74").Value = 15380422
$ws.Range("I74").Value = 22214610
$ws.Range("K74").Value = 22214610
$ws.Range("M74").Value = -22213736

$ws.Range("H77").Value = 15380422
$ws.Range("I77").Value = 22214610
$ws.Range("K77").Value = 111073050
$ws.Range("M77").Value = -111068682

$ws.Range("H116").Value = 1147.9231
$ws.Range("I116").Value = 993.9091
$ws.Range("K116").Value = 993.9091
$ws.Range("M116").Value = 1300.0909

$ws.Range("H122").Value = 1277.5
$ws.Range("I122").Value = 1341.6923
$ws.Range("K122").Value = 4025.0769
$ws.Range("M122").Value = -1575.0769

$ws.Range("H132").Value = 2199.0217
$ws.Range("I132").Value = 1404.8
$ws.Range("K132").Value = 4214.4
$ws.Range("M132").Value = -1684.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1147.9231
$ws.Range("I3").Value = 993.9091
$ws.Range("K3").Value = 993.9091
$ws.Range("M3").Value = -879.9091

$ws.Range("H7").Value = 5999.6665
$ws.Range("J7").Value = 8499.5
$ws.Range("L7").Value = 8499.5
$ws.Range("N7").Value = -8725.5

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H55").Value = 45000
$ws.Range("J55").Value = 45000
$ws.Range("L55").Value = 45000
$ws.Range("N55").Value = -45546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 54
$ws.Range("I19").Value = 54
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 54
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 116
$ws.Range("N19").ClearContents()

$ws.Range("H24").Value = 54
$ws.Range("I24").Value = 54
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 54
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 116
$ws.Range("N24").ClearContents()

$ws.Range("H109").Value = 61284
$ws.Range("J109").Value = 61284
$ws.Range("L109").Value = 61284
$ws.Range("N109").Value = -63364

$ws.Range("H132").Value = 2682.6365
$ws.Range("I132").Value = 1930.4286
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 5791.2858
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3261.2858
$ws.Range("N132").Value = -17057

$ws.Range("H134").Value = 2356.818
$ws.Range("I134").Value = 2356.818
$ws.Range("K134").Value = 7070.454000000001
$ws.Range("M134").Value = -4535.454000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 200000
$ws.Range("J37").Value = 200000
$ws.Range("L37").Value = 600000
$ws.Range("N37").Value = -600224

$ws.Range("H114").Value = 2648
$ws.Range("I114").Value = 3624.5
$ws.Range("J114").Value = 1997
$ws.Range("K114").Value = 10873.5
$ws.Range("L114").Value = 5991
$ws.Range("M114").Value = -7619.5
$ws.Range("N114").Value = -12499

$ws.Range("H129").Value = 1399.4
$ws.Range("I129").Value = 1949
$ws.Range("J129").Value = 1033
$ws.Range("K129").Value = 5847
$ws.Range("L129").Value = 3099
$ws.Range("M129").Value = -847
$ws.Range("N129").Value = -13099

$ws.Range("H131").Value = 1361.75
$ws.Range("J131").Value = 1724.25
$ws.Range("L131").Value = 5172.75
$ws.Range("N131").Value = -15252.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2255.7144
$ws.Range("I113").Value = 2058
$ws.Range("K113").Value = 2058
$ws.Range("M113").Value = 112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 466.33334
$ws.Range("I16").Value = 466.33334
$ws.Range("K16").Value = 466.33334
$ws.Range("M16").Value = -296.33334

$ws.Range("H22").Value = 967.8
$ws.Range("I22").Value = 959.5
$ws.Range("J22").Value = 973.3333
$ws.Range("K22").Value = 959.5
$ws.Range("L22").Value = 973.3333
$ws.Range("M22").Value = -664.5
$ws.Range("N22").Value = -1563.3333

$ws.Range("H27").Value = 967.8
$ws.Range("I27").Value = 959.5
$ws.Range("J27").Value = 973.3333
$ws.Range("K27").Value = 959.5
$ws.Range("L27").Value = 973.3333
$ws.Range("M27").Value = -852.5
$ws.Range("N27").Value = -1187.3333

$ws.Range("H122").Value = 4870.857
$ws.Range("I122").Value = 4019.4
$ws.Range("J122").Value = 6999.5
$ws.Range("K122").Value = 12058.2
$ws.Range("L122").Value = 20998.5
$ws.Range("M122").Value = -9608.200000000001
$ws.Range("N122").Value = -25898.5

$ws.Range("H132").Value = 4145.6924
$ws.Range("I132").Value = 3483
$ws.Range("J132").Value = 4713.7144
$ws.Range("K132").Value = 10449
$ws.Range("L132").Value = 14141.1432
$ws.Range("M132").Value = -7919
$ws.Range("N132").Value = -19201.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 80600
$ws.Range("J27").Value = 80600
$ws.Range("L27").Value = 80600
$ws.Range("N27").Value = -80738

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H132").Value = 3721.7144
$ws.Range("I132").Value = 3212.1428
$ws.Range("J132").Value = 3976.5
$ws.Range("K132").Value = 9636.428400000001
$ws.Range("L132").Value = 11929.5
$ws.Range("M132").Value = -7106.428400000001
$ws.Range("N132").Value = -16989.5
